# Swap the contents of columns C and D (codeforiati:group-name <-> codeforiati:group-code)
# for every row in the used range, including the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
